$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 0.6773940773868503
$ws.Cells.Item(2, 3).Value = 0.1887535822251891
$ws.Cells.Item(2, 4).Value = 0.02199359142910851
$ws.Cells.Item(2, 6).Value = 0.3026035079680653
$ws.Cells.Item(2, 7).Value = 0.1640119456867097
$ws.Cells.Item(2, 8).Value = 0.3333932871030996
$ws.Cells.Item(2, 14).Value = 0.8281289212184788
$ws.Cells.Item(2, 15).Value = 0.8963404678724913
$ws.Cells.Item(3, 2).Value = 0.5917434162316226
$ws.Cells.Item(3, 3).Value = 0.1717223696886379
$ws.Cells.Item(3, 4).Value = 0.01917051137701264
$ws.Cells.Item(3, 6).Value = 0.299854229014997
$ws.Cells.Item(3, 7).Value = 0.162391475602135
$ws.Cells.Item(3, 8).Value = 0.3359207621008977
$ws.Cells.Item(3, 14).Value = 0.8190346368419341
$ws.Cells.Item(3, 15).Value = 0.8979701823889457
$ws.Cells.Item(4, 2).Value = 0.5389689900325152
$ws.Cells.Item(4, 3).Value = 0.1612131111165525
$ws.Cells.Item(4, 4).Value = 0.01743127524508736
$ws.Cells.Item(4, 6).Value = 0.298434719021671
$ws.Cells.Item(4, 7).Value = 0.161589695288562
$ws.Cells.Item(4, 8).Value = 0.3376728840177776
$ws.Cells.Item(4, 14).Value = 0.8137666279988878
$ws.Cells.Item(4, 15).Value = 0.8997916614376322
$ws.Cells.Item(5, 2).Value = 0.5174182322895717
$ws.Cells.Item(5, 3).Value = 0.1569178077826621
$ws.Cells.Item(5, 4).Value = 0.01672109954305512
$ws.Cells.Item(5, 6).Value = 0.2979237117355638
$ws.Cells.Item(5, 7).Value = 0.1613113830567983
$ws.Cells.Item(5, 8).Value = 0.3384372339156272
$ws.Cells.Item(5, 14).Value = 0.8116997483698185
$ws.Cells.Item(5, 15).Value = 0.9007399732316088
$ws.Cells.Item(6, 2).Value = 0.5138370915066162
$ws.Cells.Item(6, 3).Value = 0.156203820580231
$ws.Cells.Item(6, 4).Value = 0.0166030909140531
$ws.Cells.Item(6, 6).Value = 0.2978429311933155
$ws.Cells.Item(6, 7).Value = 0.1612680892286562
$ws.Cells.Item(6, 8).Value = 0.33856719421604
$ws.Cells.Item(6, 14).Value = 0.8113613818661491
$ws.Cells.Item(6, 15).Value = 0.9009098736704146
$ws.Cells.Item(7, 2).Value = 0.5386785279663968
$ws.Cells.Item(7, 3).Value = 0.1611552340436049
$ws.Cells.Item(7, 4).Value = 0.0174217032630537
$ws.Cells.Item(7, 6).Value = 0.2984275543798489
$ws.Cells.Item(7, 7).Value = 0.1615857460358967
$ws.Cells.Item(7, 8).Value = 0.3376829884675772
$ws.Cells.Item(7, 14).Value = 0.813738429337107
$ws.Cells.Item(7, 15).Value = 0.8998036169145962
$ws.Cells.Item(8, 2).Value = 0.6479010084711092
$ws.Cells.Item(8, 3).Value = 0.1828922448564185
$ws.Cells.Item(8, 4).Value = 0.02102143695555725
$ws.Cells.Item(8, 6).Value = 0.3015997695910286
$ws.Cells.Item(8, 7).Value = 0.1634130126074211
$ws.Cells.Item(8, 8).Value = 0.3342231961404067
$ws.Cells.Item(8, 14).Value = 0.8249278988588316
$ws.Cells.Item(8, 15).Value = 0.8967318123662835
$ws.Cells.Item(9, 2).Value = 0.8605583399490229
$ws.Cells.Item(9, 3).Value = 0.2250905936805339
$ws.Cells.Item(9, 4).Value = 0.0280322226050842
$ws.Cells.Item(9, 6).Value = 0.3099556784406872
$ws.Cells.Item(9, 7).Value = 0.1685369038263147
$ws.Cells.Item(9, 8).Value = 0.3290278382903367
$ws.Cells.Item(9, 14).Value = 0.8493599595803261
$ws.Cells.Item(9, 15).Value = 0.8972384541274749
$ws.Cells.Item(10, 2).Value = 1.015796166096266
$ws.Cells.Item(10, 3).Value = 0.2558154813813189
$ws.Cells.Item(10, 4).Value = 0.03315161726075644
$ws.Cells.Item(10, 6).Value = 0.3174036466312842
$ws.Cells.Item(10, 7).Value = 0.1732519745600811
$ws.Cells.Item(10, 8).Value = 0.3261805387543788
$ws.Cells.Item(10, 14).Value = 0.8688078743485192
$ws.Cells.Item(10, 15).Value = 0.9016179939203539
$ws.Cells.Item(11, 2).Value = 1.08618640142646
$ws.Cells.Item(11, 3).Value = 0.2697291004296289
$ws.Cells.Item(11, 4).Value = 0.03547334063576102
$ws.Cells.Item(11, 6).Value = 0.3210777283609474
$ws.Cells.Item(11, 7).Value = 0.1756058347595442
$ws.Cells.Item(11, 8).Value = 0.3250959962312052
$ws.Cells.Item(11, 14).Value = 0.8779764625393796
$ws.Cells.Item(11, 15).Value = 0.9044861938670579
$ws.Cells.Item(12, 2).Value = 1.112807028855286
$ws.Cells.Item(12, 3).Value = 0.2749883682131724
$ws.Cells.Item(12, 4).Value = 0.03635144921832989
$ws.Cells.Item(12, 6).Value = 0.3225102277540799
$ws.Cells.Item(12, 7).Value = 0.1765274066515587
$ws.Cells.Item(12, 8).Value = 0.3247156211905562
$ws.Cells.Item(12, 14).Value = 0.8814942376869226
$ws.Cells.Item(12, 15).Value = 0.9056986773332198
$ws.Cells.Item(13, 2).Value = 1.107075369225413
$ws.Cells.Item(13, 3).Value = 0.273856120730926
$ws.Cells.Item(13, 4).Value = 0.03616238162229024
$ws.Cells.Item(13, 6).Value = 0.3221998795021932
$ws.Cells.Item(13, 7).Value = 0.1763275831123607
$ws.Cells.Item(13, 8).Value = 0.3247961930909682
$ws.Cells.Item(13, 14).Value = 0.880734590224705
$ws.Cells.Item(13, 15).Value = 0.9054319211655582
$ws.Cells.Item(14, 2).Value = 1.088377200785658
$ws.Cells.Item(14, 3).Value = 0.2701619766926342
$ws.Cells.Item(14, 4).Value = 0.03554560506803739
$ws.Cells.Item(14, 6).Value = 0.3211947546823524
$ws.Cells.Item(14, 7).Value = 0.1756810464953844
$ws.Cells.Item(14, 8).Value = 0.3250640948340617
$ws.Cells.Item(14, 14).Value = 0.8782649559138491
$ws.Cells.Item(14, 15).Value = 0.9045834109204236
$ws.Cells.Item(15, 2).Value = 1.076919469665995
$ws.Cells.Item(15, 3).Value = 0.2678979549421285
$ws.Cells.Item(15, 4).Value = 0.03516766965496743
$ws.Cells.Item(15, 6).Value = 0.3205844551657222
$ws.Cells.Item(15, 7).Value = 0.17528896416826
$ws.Cells.Item(15, 8).Value = 0.3252321411729966
$ws.Cells.Item(15, 14).Value = 0.8767581895739767
$ws.Cells.Item(15, 15).Value = 0.9040801413123916
$ws.Cells.Item(16, 2).Value = 1.011191251478181
$ws.Cells.Item(16, 3).Value = 0.2549048844709887
$ws.Cells.Item(16, 4).Value = 0.03299973934598199
$ws.Cells.Item(16, 6).Value = 0.3171692982268155
$ws.Cells.Item(16, 7).Value = 0.1731023633995363
$ws.Cells.Item(16, 8).Value = 0.3262556577035554
$ws.Cells.Item(16, 14).Value = 0.8682151215649583
$ws.Cells.Item(16, 15).Value = 0.901448210596854
$ws.Cells.Item(17, 2).Value = 0.9708093758545715
$ws.Cells.Item(17, 3).Value = 0.2469175501011591
$ws.Cells.Item(17, 4).Value = 0.03166792321938772
$ws.Cells.Item(17, 6).Value = 0.3151475120591343
$ws.Cells.Item(17, 7).Value = 0.1718145930409136
$ws.Cells.Item(17, 8).Value = 0.3269375323890316
$ws.Cells.Item(17, 14).Value = 0.8630563033023719
$ws.Cells.Item(17, 15).Value = 0.900058239087727
$ws.Cells.Item(18, 2).Value = 0.9475614147517604
$ws.Cells.Item(18, 3).Value = 0.2423175146695939
$ws.Cells.Item(18, 4).Value = 0.03090123138186129
$ws.Cells.Item(18, 6).Value = 0.3140115455461796
$ws.Cells.Item(18, 7).Value = 0.1710935590468168
$ws.Cells.Item(18, 8).Value = 0.3273495604826593
$ws.Cells.Item(18, 14).Value = 0.8601193860578604
$ws.Cells.Item(18, 15).Value = 0.8993411868079022
$ws.Cells.Item(19, 2).Value = 0.9396864386664561
$ws.Cells.Item(19, 3).Value = 0.2407590156271056
$ws.Cells.Item(19, 4).Value = 0.03064152990184255
$ws.Cells.Item(19, 6).Value = 0.3136315462962997
$ws.Cells.Item(19, 7).Value = 0.1708528003226277
$ws.Cells.Item(19, 8).Value = 0.3274924716907037
$ws.Cells.Item(19, 14).Value = 0.8591302130556784
$ws.Cells.Item(19, 15).Value = 0.8991125490673397
$ws.Cells.Item(20, 2).Value = 0.9751103207169649
$ws.Cells.Item(20, 3).Value = 0.2477684324434222
$ws.Cells.Item(20, 4).Value = 0.03180976668704716
$ws.Cells.Item(20, 6).Value = 0.3153599488088332
$ws.Cells.Item(20, 7).Value = 0.1719496426891496
$ws.Cells.Item(20, 8).Value = 0.3268628930146278
$ws.Cells.Item(20, 14).Value = 0.8636023355880553
$ws.Cells.Item(20, 15).Value = 0.9001976705778674
$ws.Cells.Item(21, 2).Value = 1.093870260262918
$ws.Cells.Item(21, 3).Value = 0.2712472979729625
$ws.Cells.Item(21, 4).Value = 0.0357267969759647
$ws.Cells.Item(21, 6).Value = 0.321488865279548
$ws.Cells.Item(21, 7).Value = 0.1758701285766264
$ws.Cells.Item(21, 8).Value = 0.3249845826433386
$ws.Cells.Item(21, 14).Value = 0.8789891068981888
$ws.Cells.Item(21, 15).Value = 0.9048292064440773
$ws.Cells.Item(22, 2).Value = 1.171284076883126
$ws.Cells.Item(22, 3).Value = 0.2865364241713735
$ws.Cells.Item(22, 4).Value = 0.03828049584716098
$ws.Cells.Item(22, 6).Value = 0.3257346742103451
$ws.Cells.Item(22, 7).Value = 0.1786085925073735
$ws.Cells.Item(22, 8).Value = 0.3239337209379585
$ws.Cells.Item(22, 14).Value = 0.8893121906308892
$ws.Cells.Item(22, 15).Value = 0.9085928788622084
$ws.Cells.Item(23, 2).Value = 1.129986006636955
$ws.Cells.Item(23, 3).Value = 0.2783815611016962
$ws.Cells.Item(23, 4).Value = 0.03691813381807663
$ws.Cells.Item(23, 6).Value = 0.3234466001393628
$ws.Cells.Item(23, 7).Value = 0.1771308450872695
$ws.Cells.Item(23, 8).Value = 0.3244784096036426
$ws.Cells.Item(23, 14).Value = 0.8837782793908815
$ws.Cells.Item(23, 15).Value = 0.9065165970544342
$ws.Cells.Item(24, 2).Value = 0.9731659627175304
$ws.Cells.Item(24, 3).Value = 0.2473837734567041
$ws.Cells.Item(24, 4).Value = 0.03174564240696043
$ws.Cells.Item(24, 6).Value = 0.3152638239735452
$ws.Cells.Item(24, 7).Value = 0.1718885265692407
$ws.Cells.Item(24, 8).Value = 0.3268965751502506
$ws.Cells.Item(24, 14).Value = 0.8633553841638957
$ws.Cells.Item(24, 15).Value = 0.9001343780201552
$ws.Cells.Item(25, 2).Value = 0.8031997146240428
$ws.Cells.Item(25, 3).Value = 0.2137224316142579
$ws.Cells.Item(25, 4).Value = 0.02614099429465
$ws.Cells.Item(25, 6).Value = 0.3074658435403776
$ws.Cells.Item(25, 7).Value = 0.1669846417909255
$ws.Cells.Item(25, 8).Value = 0.3302630864121028
$ws.Cells.Item(25, 14).Value = 0.8424860056958039
$ws.Cells.Item(25, 15).Value = 0.896399542589009
